# Updated symbol list on Mon Feb 13 07:45:24 UTC 2023 with GitHub Actions
#
# Applies the cryptocurrency price/volume refresh + row-shift described by
# the commit: several coins (FTXToken, MXToken, LiechtensteinCryptoassets-
# Exchange, WazirX, MandalaExchangeToken, BitrueCoin, BitMartToken,
# BitForexToken, TigerCash, LEO, GateToken) moved down one rank (their
# Coin name + Link moved to the next row), and most Price / Volume(1h)
# cells across the sheet were refreshed with newer scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell, preserving the original
# inlineStr/text cell semantics even for values that look numeric
# (prices like "311.53" or percentages like "1.31%"). Excel's COM layer
# auto-coerces plain numeric-looking .Value assignments into real numbers,
# so for those we briefly force a Text number format, assign the value,
# then restore the cell style to Normal (removing the temporary format)
# so no stray formatting is left behind on the cell.
function Set-CellText($sheet, $addr, $text) {
    $range = $sheet.Range($addr)
    if ($text -match '^-?[0-9]') {
        $range.NumberFormat = "@"
        $range.Value = $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

Set-CellText $ws "D2" "311.53"
Set-CellText $ws "E2" "1.31%"
Set-CellText $ws "D3" "41.00"
Set-CellText $ws "E3" "-0.06%"
Set-CellText $ws "D4" "5.164"
Set-CellText $ws "E4" "-1.23%"
Set-CellText $ws "E5" "-0.72%"
Set-CellText $ws "B6" "GateToken"
Set-CellText $ws "C6" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-CellText $ws "D6" "4.322"
Set-CellText $ws "E6" "0.27%"
Set-CellText $ws "B7" "FTXToken"
Set-CellText $ws "C7" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-CellText $ws "D7" "1.681"
Set-CellText $ws "E7" "2.42%"
Set-CellText $ws "B8" "MXToken"
Set-CellText $ws "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-CellText $ws "D8" "0.9322"
Set-CellText $ws "E8" "1.90%"
Set-CellText $ws "B9" "LiechtensteinCryptoassetsExchange"
Set-CellText $ws "C9" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-CellText $ws "D9" "0.1200"
Set-CellText $ws "E9" "-3.25%"
Set-CellText $ws "B10" "WazirX"
Set-CellText $ws "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-CellText $ws "D10" "0.1821"
Set-CellText $ws "E10" "-0.22%"
Set-CellText $ws "B11" "MandalaExchangeToken"
Set-CellText $ws "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-CellText $ws "D11" "0.09012"
Set-CellText $ws "E11" "-1.98%"
Set-CellText $ws "B12" "BitrueCoin"
Set-CellText $ws "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-CellText $ws "D12" "0.04138"
Set-CellText $ws "E12" "0.75%"
Set-CellText $ws "B13" "BitMartToken"
Set-CellText $ws "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-CellText $ws "D13" "0.1055"
Set-CellText $ws "E13" "0.31%"
Set-CellText $ws "B14" "BitForexToken"
Set-CellText $ws "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-CellText $ws "D14" "0.001282"
Set-CellText $ws "E14" "1.98%"
Set-CellText $ws "B15" "TigerCash"
Set-CellText $ws "C15" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-CellText $ws "D15" "0.005860"
Set-CellText $ws "E15" "0.05%"
Set-CellText $ws "B16" "LEO"
Set-CellText $ws "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-CellText $ws "D16" "3.337"
Set-CellText $ws "E16" "-0.28%"
Set-CellText $ws "E17" "-0.49%"
Set-CellText $ws "D18" "0.3350"
Set-CellText $ws "E18" "0.42%"
Set-CellText $ws "D19" "7.584"
Set-CellText $ws "E19" "1.71%"
Set-CellText $ws "D20" "0.1340"
Set-CellText $ws "E20" "-3.82%"
Set-CellText $ws "D21" "0.2836"
Set-CellText $ws "E21" "4.43%"
Set-CellText $ws "D22" "0.03986"
Set-CellText $ws "E22" "-1.42%"
Set-CellText $ws "D23" "0.001280"
Set-CellText $ws "E23" "1.22%"
Set-CellText $ws "D24" "0.004065"
Set-CellText $ws "E24" "-5.46%"
Set-CellText $ws "D25" "0.0001350"
Set-CellText $ws "E25" "6.07%"
Set-CellText $ws "D26" "0.0003040"
Set-CellText $ws "D38" "0.02436"
Set-CellText $ws "E38" "-1.77%"
Set-CellText $ws "D39" "0.05182"
Set-CellText $ws "E39" "-2.86%"
Set-CellText $ws "D40" "0.007705"
Set-CellText $ws "E40" "-1.90%"
Set-CellText $ws "E41" "-0.76%"
Set-CellText $ws "D42" "0.007591"
Set-CellText $ws "E42" "15.21%"
Set-CellText $ws "E43" "72.41%"
Set-CellText $ws "D44" "0.008478"
Set-CellText $ws "E44" "10.60%"
Set-CellText $ws "D45" "0.3379"
Set-CellText $ws "E45" "0.96%"
Set-CellText $ws "E46" "-2.03%"
Set-CellText $ws "E47" "-0.21%"
Set-CellText $ws "D48" "0.2701"
Set-CellText $ws "E48" "-27.09%"
Set-CellText $ws "D49" "0.004201"
Set-CellText $ws "E49" "35.20%"
Set-CellText $ws "D50" "0.00002100"
Set-CellText $ws "E50" "-0.21%"
Set-CellText $ws "D51" "0.0002000"
Set-CellText $ws "E51" "-0.21%"

Write-Output "Applied 93 cell updates to Sheet1"
